$wb = $excel.ActiveWorkbook

# ---- Sheet: LP1912 ----
$ws = $wb.Worksheets.Item('LP1912')
$ws.Range("A2").Value = 'Última actualización: 08:56:26'
$ws.Range("A3").Value = 'Total filas: 134'

$rows_1 = @(
    @(6, '01:10:32', '01:12', '215_ALUAR', 2, 'LP1912'),
    @(7, '01:55:51', '01:58', '14_ABASTO', 3, 'LP1912'),
    @(8, '02:48:52', '02:57', '215_ALUAR', 9, 'LP1912'),
    @(9, '02:58:51', '02:58', '215_ALUAR', 0, 'LP1912'),
    @(10, '01:55:51', '03:12', '215_ALUAR', 77, 'LP1912'),
    @(11, '03:35:49', '03:38', '14_ABASTO', 3, 'LP1912'),
    @(12, '02:58:51', '03:48', '14_ABASTO', 50, 'LP1912'),
    @(13, '02:21:47', '03:56', '14_ABASTO', 95, 'LP1912'),
    @(14, '04:01:13', '04:01', '81_EL PELIGRO', 0, 'LP1912'),
    @(15, '03:35:49', '04:45', '215A_EL PATO', 70, 'LP1912'),
    @(16, '04:35:25', '04:46', '215A_EL PATO', 11, 'LP1912'),
    @(17, '04:48:57', '04:53', '11_ETCHEVERRY', 5, 'LP1912'),
    @(18, '04:48:57', '05:14', '14_ABASTO', 26, 'LP1912'),
    @(19, '04:35:25', '05:15', '14_ABASTO', 40, 'LP1912'),
    @(20, '04:56:11', '05:16', '17_ROMERO', 20, 'LP1912'),
    @(21, '04:48:57', '05:21', '23_HERNANDEZ', 33, 'LP1912'),
    @(22, '04:56:11', '05:22', '23_HERNANDEZ', 26, 'LP1912'),
    @(23, '05:21:16', '05:26', '23_HERNANDEZ', 5, 'LP1912'),
    @(24, '03:35:49', '05:28', '14_ABASTO', 113, 'LP1912'),
    @(25, '04:48:57', '05:34', '215B_EL PATO', 46, 'LP1912'),
    @(26, '05:21:16', '05:35', '215B_EL PATO', 14, 'LP1912'),
    @(27, '04:01:13', '05:37', '14_ABASTO', 96, 'LP1912'),
    @(28, '05:21:16', '05:46', '15_ABASTO', 25, 'LP1912'),
    @(29, '05:52:07', '05:54', '10_OLMOS', 2, 'LP1912'),
    @(30, '05:21:16', '06:04', '16_SANTA ANA', 43, 'LP1912'),
    @(31, '05:21:16', '06:11', '215A_EL PATO', 50, 'LP1912'),
    @(32, '05:52:07', '06:12', '215A_EL PATO', 20, 'LP1912'),
    @(33, '04:48:57', '06:13', '225_HARAS DEL SUR', 85, 'LP1912'),
    @(34, '05:52:07', '06:14', '225_HARAS DEL SUR', 22, 'LP1912'),
    @(35, '04:48:57', '06:20', '26_HERNANDEZ', 92, 'LP1912'),
    @(36, '06:21:22', '06:21', '26_HERNANDEZ', 0, 'LP1912'),
    @(37, '04:48:57', '06:26', '23_HERNANDEZ', 98, 'LP1912'),
    @(38, '05:52:07', '06:27', '23_HERNANDEZ', 35, 'LP1912'),
    @(39, '06:21:22', '06:29', '23_HERNANDEZ', 8, 'LP1912'),
    @(40, '06:21:22', '06:29', '86_EST CHICA-ESC AGRARIA', 8, 'LP1912'),
    @(41, '05:52:07', '06:30', '86_EST CHICA-ESC AGRARIA', 38, 'LP1912'),
    @(42, '06:21:22', '06:31', '16_SANTA ANA', 10, 'LP1912'),
    @(43, '04:48:57', '06:43', '225_C ROCA-H SUR', 115, 'LP1912'),
    @(44, '06:21:22', '06:44', '225_C ROCA-H SUR', 23, 'LP1912'),
    @(45, '06:21:22', '06:46', '215C_EL PATO', 25, 'LP1912'),
    @(46, '05:52:07', '06:47', '215C_EL PATO', 55, 'LP1912'),
    @(47, '06:59:37', '06:59', '23_HERNANDEZ', 0, 'LP1912'),
    @(48, '06:59:37', '06:59', '14_ABASTO', 0, 'LP1912'),
    @(49, '05:52:07', '07:00', '14_ABASTO', 68, 'LP1912'),
    @(50, '06:49:33', '07:01', '16_SANTA ANA', 12, 'LP1912'),
    @(51, '06:49:33', '07:04', '23_HERNANDEZ', 15, 'LP1912'),
    @(52, '06:59:37', '07:05', '15_ABASTO', 6, 'LP1912'),
    @(53, '05:52:07', '07:05', '23_HERNANDEZ', 73, 'LP1912'),
    @(54, '06:59:37', '07:07', '225_GOMEZ', 8, 'LP1912'),
    @(55, '06:59:37', '07:11', '215A_EL PATO', 12, 'LP1912'),
    @(56, '05:52:07', '07:12', '215A_EL PATO', 80, 'LP1912'),
    @(57, '06:59:37', '07:15', '11_ETCHEVERRY', 16, 'LP1912'),
    @(58, '06:59:37', '07:16', '16_SANTA ANA', 17, 'LP1912'),
    @(59, '05:52:07', '07:16', '11_ETCHEVERRY', 84, 'LP1912'),
    @(60, '06:59:37', '07:21', '26_HERNANDEZ', 22, 'LP1912'),
    @(61, '06:59:37', '07:23', '10_OLMOS', 24, 'LP1912'),
    @(62, '07:28:14', '07:30', '11_ETCHEVERRY', 2, 'LP1912'),
    @(63, '06:59:37', '07:31', '11_ETCHEVERRY', 32, 'LP1912'),
    @(64, '06:59:37', '07:31', '16_SANTA ANA', 32, 'LP1912'),
    @(65, '05:52:07', '07:32', '16_SANTA ANA', 100, 'LP1912'),
    @(66, '07:28:14', '07:32', '84_COLONIA URQUIZA-ESC 49', 4, 'LP1912'),
    @(67, '05:52:07', '07:32', '11_ETCHEVERRY', 100, 'LP1912'),
    @(68, '07:28:14', '07:35', '23_HERNANDEZ', 7, 'LP1912'),
    @(69, '06:59:37', '07:36', '27_EL RETIRO', 37, 'LP1912'),
    @(70, '07:28:14', '07:37', '27_EL RETIRO', 9, 'LP1912'),
    @(71, '07:28:14', '07:39', '10_OLMOS', 11, 'LP1912'),
    @(72, '06:59:37', '07:47', '14_ABASTO', 48, 'LP1912'),
    @(73, '07:28:14', '07:47', '16_SANTA ANA', 19, 'LP1912'),
    @(74, '07:28:14', '07:48', '14_ABASTO', 20, 'LP1912'),
    @(75, '07:51:34', '07:51', '215D_EL PATO', 0, 'LP1912'),
    @(76, '07:51:34', '07:51', '10_OLMOS', 0, 'LP1912'),
    @(77, '07:28:14', '07:55', '10_OLMOS', 27, 'LP1912'),
    @(78, '07:28:14', '08:00', '23_HERNANDEZ', 32, 'LP1912'),
    @(79, '07:51:34', '08:01', '23_HERNANDEZ', 10, 'LP1912'),
    @(80, '07:51:34', '08:03', '11_ETCHEVERRY', 12, 'LP1912'),
    @(81, '06:59:37', '08:06', '23_HERNANDEZ', 67, 'LP1912'),
    @(82, '07:51:34', '08:10', '16_SANTA ANA', 19, 'LP1912'),
    @(83, '07:28:14', '08:11', '16_SANTA ANA', 43, 'LP1912'),
    @(84, '07:51:34', '08:12', '15_ABASTO', 21, 'LP1912'),
    @(85, '07:51:34', '08:13', '10_OLMOS', 22, 'LP1912'),
    @(86, '08:13:38', '08:21', '26_HERNANDEZ', 8, 'LP1912'),
    @(87, '08:13:38', '08:22', '16_P MOR-SANTA ANA', 9, 'LP1912'),
    @(88, '07:28:14', '08:23', '16_P MOR-SANTA ANA', 55, 'LP1912'),
    @(89, '08:13:38', '08:23', '215B_EL PATO', 10, 'LP1912'),
    @(90, '08:13:38', '08:27', '84_COLONIA URQUIZA-ESC 49', 14, 'LP1912'),
    @(91, '07:51:34', '08:30', '23_HERNANDEZ', 39, 'LP1912'),
    @(92, '08:13:38', '08:33', '10_OLMOS', 20, 'LP1912'),
    @(93, '08:13:38', '08:36', '23_HERNANDEZ', 23, 'LP1912'),
    @(94, '08:34:05', '08:42', '81_EL PELIGRO', 8, 'LP1912'),
    @(95, '08:13:38', '08:43', '14_ABASTO', 30, 'LP1912'),
    @(96, '08:34:05', '08:44', '14_ABASTO', 10, 'LP1912'),
    @(97, '08:48:01', '08:53', '10_OLMOS', 5, 'LP1912'),
    @(98, '08:34:05', '08:54', '17_ROMERO', 20, 'LP1912'),
    @(99, '08:13:38', '09:01', '23_HERNANDEZ', 48, 'LP1912'),
    @(100, '08:56:26', '09:01', '215A_EL PATO', 5, 'LP1912'),
    @(101, '08:48:01', '09:02', '215A_EL PATO', 14, 'LP1912'),
    @(102, '08:56:26', '09:03', '11_ETCHEVERRY', 7, 'LP1912'),
    @(103, '08:56:26', '09:04', '23_HERNANDEZ', 8, 'LP1912'),
    @(104, '08:48:01', '09:04', '11_ETCHEVERRY', 16, 'LP1912'),
    @(105, '08:48:01', '09:05', '23_HERNANDEZ', 17, 'LP1912'),
    @(106, '08:56:26', '09:10', '16_P MOR-SANTA ANA', 14, 'LP1912'),
    @(107, '08:48:01', '09:11', '16_P MOR-SANTA ANA', 23, 'LP1912'),
    @(108, '08:56:26', '09:13', '10_OLMOS', 17, 'LP1912'),
    @(109, '08:56:26', '09:16', '27_EL RETIRO', 20, 'LP1912'),
    @(110, '08:48:01', '09:17', '27_EL RETIRO', 29, 'LP1912'),
    @(111, '08:56:26', '09:21', '26_HERNANDEZ', 25, 'LP1912'),
    @(112, '08:13:38', '09:22', '17_ROMERO', 69, 'LP1912'),
    @(113, '07:28:14', '09:23', '17_ROMERO', 115, 'LP1912'),
    @(114, '08:56:26', '09:23', '11_ETCHEVERRY', 27, 'LP1912'),
    @(115, '08:56:26', '09:23', '16_SANTA ANA', 27, 'LP1912'),
    @(116, '08:48:01', '09:24', '11_ETCHEVERRY', 36, 'LP1912'),
    @(117, '08:56:26', '09:32', '15_ABASTO', 36, 'LP1912'),
    @(118, '08:56:26', '09:33', '10_OLMOS', 37, 'LP1912'),
    @(119, '08:56:26', '09:34', '23_HERNANDEZ', 38, 'LP1912'),
    @(120, '08:56:26', '09:34', '16_SANTA ANA', 38, 'LP1912'),
    @(121, '08:48:01', '09:35', '23_HERNANDEZ', 47, 'LP1912'),
    @(122, '08:48:01', '09:35', '16_SANTA ANA', 47, 'LP1912'),
    @(123, '08:13:38', '09:41', '215C_EL PATO', 88, 'LP1912'),
    @(124, '08:56:26', '09:42', '215C_EL PATO', 46, 'LP1912'),
    @(125, '08:56:26', '09:43', '14_ABASTO', 47, 'LP1912'),
    @(126, '08:48:01', '09:44', '14_ABASTO', 56, 'LP1912'),
    @(127, '08:56:26', '09:52', '15_ABASTO', 56, 'LP1912'),
    @(128, '08:56:26', '09:53', '10_OLMOS', 57, 'LP1912'),
    @(129, '08:13:38', '09:58', '16_SANTA ANA', 105, 'LP1912'),
    @(130, '08:56:26', '10:10', '16_P MOR-SANTA ANA', 74, 'LP1912'),
    @(131, '08:48:01', '10:11', '16_P MOR-SANTA ANA', 83, 'LP1912'),
    @(132, '08:34:05', '10:12', '15_ABASTO', 98, 'LP1912'),
    @(133, '08:56:26', '10:21', '26_HERNANDEZ', 85, 'LP1912'),
    @(134, '08:56:26', '10:22', '17_ROMERO', 86, 'LP1912'),
    @(135, '08:56:26', '10:26', '215A_EL PATO', 90, 'LP1912'),
    @(136, '08:48:01', '10:27', '215A_EL PATO', 99, 'LP1912'),
    @(137, '08:56:26', '10:42', '17_ROMERO', 106, 'LP1912'),
    @(138, '08:56:26', '10:43', '14_ABASTO', 107, 'LP1912'),
    @(139, '08:48:01', '10:44', '14_ABASTO', 116, 'LP1912'),
)

foreach ($row in $rows_1) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
}

# ---- Sheet: LP1912-215 ----
$ws = $wb.Worksheets.Item('LP1912-215')
$ws.Range("A2").Value = 'Última actualización: 08:56:26'
$ws.Range("A3").Value = 'Total filas: 22'

$rows_2 = @(
    @(6, '01:10:32', '01:12', '215_ALUAR', 2, 'LP1912'),
    @(7, '02:48:52', '02:57', '215_ALUAR', 9, 'LP1912'),
    @(8, '02:58:51', '02:58', '215_ALUAR', 0, 'LP1912'),
    @(9, '01:55:51', '03:12', '215_ALUAR', 77, 'LP1912'),
    @(10, '03:35:49', '04:45', '215A_EL PATO', 70, 'LP1912'),
    @(11, '04:35:25', '04:46', '215A_EL PATO', 11, 'LP1912'),
    @(12, '04:48:57', '05:34', '215B_EL PATO', 46, 'LP1912'),
    @(13, '05:21:16', '05:35', '215B_EL PATO', 14, 'LP1912'),
    @(14, '05:21:16', '06:11', '215A_EL PATO', 50, 'LP1912'),
    @(15, '05:52:07', '06:12', '215A_EL PATO', 20, 'LP1912'),
    @(16, '06:21:22', '06:46', '215C_EL PATO', 25, 'LP1912'),
    @(17, '05:52:07', '06:47', '215C_EL PATO', 55, 'LP1912'),
    @(18, '06:59:37', '07:11', '215A_EL PATO', 12, 'LP1912'),
    @(19, '05:52:07', '07:12', '215A_EL PATO', 80, 'LP1912'),
    @(20, '07:51:34', '07:51', '215D_EL PATO', 0, 'LP1912'),
    @(21, '08:13:38', '08:23', '215B_EL PATO', 10, 'LP1912'),
    @(22, '08:56:26', '09:01', '215A_EL PATO', 5, 'LP1912'),
    @(23, '08:48:01', '09:02', '215A_EL PATO', 14, 'LP1912'),
    @(24, '08:13:38', '09:41', '215C_EL PATO', 88, 'LP1912'),
    @(25, '08:56:26', '09:42', '215C_EL PATO', 46, 'LP1912'),
    @(26, '08:56:26', '10:26', '215A_EL PATO', 90, 'LP1912'),
    @(27, '08:48:01', '10:27', '215A_EL PATO', 99, 'LP1912'),
)

foreach ($row in $rows_2) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
}

# ---- Sheet: 6203-6173 ----
$ws = $wb.Worksheets.Item('6203-6173')
$ws.Range("A2").Value = 'Última actualización: 08:56:26'
$ws.Range("A3").Value = 'Total filas: 26'

$rows_3 = @(
    @(6, '00:07:17', '00:08', '215A_LA PLATA', 1, 'L6173'),
    @(7, '04:48:57', '05:43', '215A_LA PLATA', 55, 'L6173'),
    @(8, '05:21:16', '05:44', '215A_LA PLATA', 23, 'L6173'),
    @(9, '04:48:57', '06:08', '215A_LA PLATA', 80, 'L6173'),
    @(10, '05:21:16', '06:09', '215A_LA PLATA', 48, 'L6173'),
    @(11, '05:52:07', '06:13', '215A_LA PLATA', 21, 'L6173'),
    @(12, '04:48:57', '06:32', '215C_LA PLATA', 104, 'L6203'),
    @(13, '06:21:22', '06:33', '215C_LA PLATA', 12, 'L6203'),
    @(14, '06:59:37', '06:59', '215B_LP-P MOR-1 Y 57', 0, 'L6173'),
    @(15, '06:49:33', '07:00', '215B_LP-P MOR-1 Y 57', 11, 'L6173'),
    @(16, '07:28:14', '07:35', '215A_LA PLATA', 7, 'L6173'),
    @(17, '06:59:37', '08:06', '215C_LA PLATA', 67, 'L6203'),
    @(18, '06:49:33', '08:07', '215C_LA PLATA', 78, 'L6203'),
    @(19, '07:28:14', '08:10', '215C_LA PLATA', 42, 'L6203'),
    @(20, '07:51:34', '08:11', '215C_LA PLATA', 20, 'L6203'),
    @(21, '08:13:38', '08:16', '215C_LA PLATA', 3, 'L6203'),
    @(22, '06:49:33', '08:33', '215A_LA PLATA', 104, 'L6173'),
    @(23, '07:28:14', '08:38', '215A_LA PLATA', 70, 'L6173'),
    @(24, '07:51:34', '08:40', '215A_LA PLATA', 49, 'L6173'),
    @(25, '08:13:38', '08:45', '215A_LA PLATA', 32, 'L6173'),
    @(26, '08:34:05', '08:46', '215A_LA PLATA', 12, 'L6173'),
    @(27, '08:13:38', '09:08', '215D_LA PLATA', 55, 'L6203'),
    @(28, '08:56:26', '09:09', '215D_LA PLATA', 13, 'L6203'),
    @(29, '08:48:01', '09:10', '215D_LA PLATA', 22, 'L6203'),
    @(30, '08:56:26', '10:03', '215B_LP-P MOR-40 Y 115', 67, 'L6173'),
    @(31, '08:56:26', '10:54', '215A_LA PLATA', 118, 'L6173'),
)

foreach ($row in $rows_3) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
}
